$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 3 becomes the "no CO2 price" (all-zero) scenario row that used to live in row 2's
# sibling; update id_scenario / id_sector and zero out the yearly price columns.
$ws1.Range("A3").Value = 1
$ws1.Range("C3").Value = 6
$ws1.Range("E3:AJ3").Value = 0

# Row 4 (the old third scenario row) is no longer needed - remove it entirely.
$ws1.Rows.Item(4).Delete()

# Remove the now-unused lookup sheet (its comment/scenario-name strings go with it).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()

# Update the active selection on Sheet1 to match the edited workbook's saved view.
$ws1.Activate()
$ws1.Range("F11").Select()
